# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2304"
#   "<name>_new" -> "<name>_FV2310"
# Then expose the data range as an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row cells -------------------------------------------------
$headerRow = 1
$lastCol = 21  # columns A..U

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item($headerRow, $col)
    $current = [string]$cell.Value2
    if ($current -ne $null) {
        if ($current.EndsWith("_old")) {
            $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2304"
        } elseif ($current.EndsWith("_new")) {
            $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2310"
        }
    }
}

# --- 2. Turn the data range into a native Excel Table ---------------------------
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row -----------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
